# edit.ps1
# Applies the "Covered more COR aspects and completed section." change:
#   1. Adds three new bullet paragraphs after the Chain-of-Responsibility
#      intro paragraph (with the "_GoBack" bookmark relocated into the
#      third of them, mid-run, exactly as in the target XML).
#   2. Moves <w:lastRenderedPageBreak/> from the "Template Method" heading
#      run to the run right before it (the "Strategy" section's ">" run).

$d = $word.ActiveDocument

function Find-ParaIndex($doc, $searchText, $startFrom) {
    for ($i = $startFrom; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $searchText) {
            return $i
        }
    }
    return -1
}

function Escape-Xml($s) {
    return ($s -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;')
}

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraPPr = '<w:pPr><w:pStyle w:val="CNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr>'

function Set-ParagraphRuns($doc, $paraIndex, [string[]]$runTexts) {
    # Fills the (already-created, empty) paragraph at $paraIndex with one
    # <w:r><w:t xml:space="preserve">...</w:t></w:r> per entry of $runTexts,
    # keeping the runs distinct (not merged into one run) and preserving
    # the paragraph's own pPr.
    $runsXml = ""
    foreach ($t in $runTexts) {
        $runsXml += "<w:r><w:t xml:space=`"preserve`">" + (Escape-Xml $t) + "</w:t></w:r>"
    }
    $xml = $pkgHeader + "<w:p>" + $paraPPr + $runsXml + "</w:p>" + $pkgFooter
    $doc.Paragraphs.Item($paraIndex).Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Locate the paragraph that ends the Chain-of-Responsibility intro.
# ---------------------------------------------------------------------
$introText = "The chain of responsibility design pattern is a chain of components " + `
    "who all get a chance to process a command or query, optionally having a " + `
    "default processing implementation and an ability to terminate the processing chain."
$introIdx = Find-ParaIndex $d $introText 1
if ($introIdx -eq -1) { throw "Could not find the Chain of Responsibility intro paragraph" }

# Remove the hidden _GoBack bookmark sitting at the end of that paragraph;
# it will be recreated (same id, since id 0 is freed) inside the new third
# paragraph below.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2) Insert three new empty bullet paragraphs right after the intro
#    paragraph. InsertParagraphAfter() clones the source paragraph's pPr
#    (CNormal / numPr ilvl=0,numId=3 / spacing after=0), which is exactly
#    the formatting the new paragraphs need.
# ---------------------------------------------------------------------
$introRange = $d.Paragraphs.Item($introIdx).Range
$introRange.InsertParagraphAfter()
$introRange.InsertParagraphAfter()
$introRange.InsertParagraphAfter()

$p1 = $introIdx + 1
$p2 = $introIdx + 2
$p3 = $introIdx + 3

# ---------------------------------------------------------------------
# 3) Fill in the first two new paragraphs (plain runs, no bookmark).
# ---------------------------------------------------------------------
Set-ParagraphRuns $d $p1 @(
    "Example: You click a button on a form. The button can handle the event, or can pass it onto its parent ",
    "- ",
    "the group box. The group box can handle the event, or pass it onto its parent ",
    "- ",
    "the window."
)

Set-ParagraphRuns $d $p2 @(
    "Example: A game creature may have multiple boosters that affect its stats."
)

# ---------------------------------------------------------------------
# 4) Fill in the third new paragraph, re-inserting the _GoBack bookmark
#    between the third and fourth runs, exactly as in the target XML.
# ---------------------------------------------------------------------
$run1 = Escape-Xml "To use the pattern, "
$run2 = Escape-Xml "store a list "
$run3 = Escape-Xml ("of objects and call their handle method when the handle method is called in the COR object. " + `
    "Call each object" + [char]0x2019 + "s method one by one and supply it with information relevant to the event. " + `
    "If the method indicates that it doesn" + [char]0x2019 + "t want to propagate")
$run4 = Escape-Xml " the event further up the chain, then stop."

$p3Xml = $pkgHeader + "<w:p>" + $paraPPr + `
    "<w:r><w:t xml:space=`"preserve`">$run1</w:t></w:r>" + `
    "<w:r><w:t xml:space=`"preserve`">$run2</w:t></w:r>" + `
    "<w:r><w:t>$run3</w:t></w:r>" + `
    "<w:bookmarkStart w:id=`"99`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"99`"/>" + `
    "<w:r><w:t xml:space=`"preserve`">$run4</w:t></w:r>" + `
    "</w:p>" + $pkgFooter

$d.Paragraphs.Item($p3).Range.InsertXML($p3Xml)

# ---------------------------------------------------------------------
# 5) Move <w:lastRenderedPageBreak/> from the "Template Method" heading
#    run to the run right before it (the ">" run under "Strategy").
# ---------------------------------------------------------------------
$stratHeadingIdx = Find-ParaIndex $d "Strategy" 1
if ($stratHeadingIdx -eq -1) { throw "Could not find the Strategy heading paragraph" }
$stratBodyIdx = $stratHeadingIdx + 1   # the ">" bullet right under "Strategy"

$stratBodyFull = $d.Paragraphs.Item($stratBodyIdx).Range
$stratBodyText = $d.Range($stratBodyFull.Start, $stratBodyFull.End - 1)
$stratXml = $pkgHeader + '<w:p><w:r><w:lastRenderedPageBreak/><w:t>&gt;</w:t></w:r></w:p>' + $pkgFooter
$stratBodyText.InsertXML($stratXml)

$tmHeadingIdx = Find-ParaIndex $d "Template Method" 1
if ($tmHeadingIdx -eq -1) { throw "Could not find the Template Method heading paragraph" }
$tmFull = $d.Paragraphs.Item($tmHeadingIdx).Range
$tmText = $d.Range($tmFull.Start, $tmFull.End - 1)
$tmXml = $pkgHeader + '<w:p><w:r><w:t>Template Method</w:t></w:r></w:p>' + $pkgFooter
$tmText.InsertXML($tmXml)

Write-Output "Edit complete."
